# "Added and Updated Czech Test Data"
#
# - Belgium sheet: no longer the active/selected tab; selection becomes the
#   full used range (A1:D11) instead of the old single-cell B6 selection.
# - A new "Czech" sheet is appended after "Belgium" (duplicated from the
#   Belgium sheet, same layout/styles/merged cells) and becomes the active
#   tab, with its own market name + ticket-id strings and the cursor left on
#   B4.

$wb = $excel.ActiveWorkbook

# --- Belgium: drop the old single-cell selection, select the whole table ---
$belgium = $wb.Worksheets.Item("Belgium")
$belgium.Select() | Out-Null
$belgium.Range("A1:D11").Select() | Out-Null

# --- Czech: duplicate Belgium's layout/styles, place it after Belgium ---
$belgium.Copy([System.Reflection.Missing]::Value, $belgium)
$czech = $wb.Worksheets.Item($wb.Worksheets.Count)
$czech.Name = "Czech"

# New market-specific content (new shared strings "Czech Market" / ticket id)
$czech.Range("B2").Value = "Czech Market"
$czech.Range("B4").Value = "NGC-3477/T1734"

# Czech's columns B:D are narrower than Belgium's (column A stays the same)
$czech.Columns.Item(2).ColumnWidth = 27.7
$czech.Columns.Item(3).ColumnWidth = 12.5
$czech.Columns.Item(4).ColumnWidth = 18.5

# Czech becomes the active sheet/tab with the cursor on B4
$czech.Select() | Out-Null
$czech.Range("B4").Select() | Out-Null
